$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D values look numeric but must remain stored as text (as in the original).
# Force text format before assignment, then restore Normal style so no stray
# number-format is left applied to the cell.

$cell = $ws.Range("D2")
$cell.NumberFormat = "@"
$cell.Value = "51.564.18"
$cell.Style = "Normal"
$ws.Range("E2").Value = "  -0.58%  "

$cell = $ws.Range("D3")
$cell.NumberFormat = "@"
$cell.Value = "2.779.47"
$cell.Style = "Normal"
$ws.Range("E3").Value = "  -0.09%  "

$ws.Range("E4").Value = "  +0.01%  "

$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = "352.18"
$cell.Style = "Normal"
$ws.Range("E5").Value = "  -1.73%  "

$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = "108.44"
$cell.Style = "Normal"
$ws.Range("E6").Value = "  -0.76%  "

$cell = $ws.Range("D7")
$cell.NumberFormat = "@"
$cell.Value = "0.550"
$cell.Style = "Normal"

$cell = $ws.Range("D8")
$cell.NumberFormat = "@"
$cell.Value = "0.999"
$cell.Style = "Normal"
$ws.Range("E8").Value = "  -0.02%  "

$cell = $ws.Range("D9")
$cell.NumberFormat = "@"
$cell.Value = "0.595"
$cell.Style = "Normal"
$ws.Range("E9").Value = "  -0.02%  "

$cell = $ws.Range("D10")
$cell.NumberFormat = "@"
$cell.Value = "39.71"
$cell.Style = "Normal"
$ws.Range("E10").Value = "  -0.95%  "

$ws.Range("E12").Value = "  +2.97%  "

$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = "7.66"
$cell.Style = "Normal"
$ws.Range("E14").Value = "  +0.88%  "

$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = "3.214.60"
$cell.Style = "Normal"
$ws.Range("E15").Value = "  -0.08%  "

$cell = $ws.Range("D16")
$cell.NumberFormat = "@"
$cell.Value = "2.765.16"
$cell.Style = "Normal"
$ws.Range("E16").Value = "  -1.23%  "

$cell = $ws.Range("D18")
$cell.NumberFormat = "@"
$cell.Value = "51.528.09"
$cell.Style = "Normal"

$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = "7.66"
$cell.Style = "Normal"
$ws.Range("E19").Value = "  +3.14%  "

$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = "3.11"
$cell.Style = "Normal"
$ws.Range("E20").Value = "  -0.58%  "

$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value = "13.13"
$cell.Style = "Normal"
$ws.Range("E21").Value = "  +0.72%  "

$ws.Range("E22").Value = "  -1.65%  "

$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = "69.93"
$cell.Style = "Normal"
$ws.Range("E23").Value = "  -0.41%  "

$cell = $ws.Range("D24")
$cell.NumberFormat = "@"
$cell.Value = "266.64"
$cell.Style = "Normal"
$ws.Range("E24").Value = "  -2.75%  "

$ws.Range("E25").Value = "  -1.02%  "

$ws.Range("E26").Value = "  +0.06%  "

$cell = $ws.Range("D27")
$cell.NumberFormat = "@"
$cell.Value = "26.10"
$cell.Style = "Normal"
$ws.Range("E27").Value = "  -2.20%  "

$ws.Range("E28").Value = "  +12.23%  "

$ws.Range("E29").Value = "  +0.59%  "

$ws.Range("E30").Value = "  -2.45%  "

$cell = $ws.Range("D31")
$cell.NumberFormat = "@"
$cell.Value = "36.28"
$cell.Style = "Normal"
$ws.Range("E31").Value = "  +6.99%  "

$cell = $ws.Range("D32")
$cell.NumberFormat = "@"
$cell.Value = "6.22"
$cell.Style = "Normal"
$ws.Range("E32").Value = "  +9.25%  "

$cell = $ws.Range("D33")
$cell.NumberFormat = "@"
$cell.Value = "51.94"
$cell.Style = "Normal"
$ws.Range("E33").Value = "  +0.99%  "

$ws.Range("E34").Value = "  -2.71%  "

$ws.Range("E35").Value = "  +5.47%  "

$ws.Range("E36").Value = "  -1.88%  "

$ws.Range("E37").Value = "  -0.03%  "

$cell = $ws.Range("D38")
$cell.NumberFormat = "@"
$cell.Value = "18.48"
$cell.Style = "Normal"
$ws.Range("E38").Value = "  +2.13%  "

$ws.Range("E39").Value = "  -2.62%  "

$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = "1.97"
$cell.Style = "Normal"
$ws.Range("E40").Value = "  -1.60%  "

$cell = $ws.Range("D41")
$cell.NumberFormat = "@"
$cell.Value = "2.55"
$cell.Style = "Normal"
$ws.Range("E41").Value = "  -0.18%  "

$ws.Range("E42").Value = "  -0.85%  "

$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = "120.46"
$cell.Style = "Normal"
$ws.Range("E43").Value = "  -1.11%  "

$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = "22.00"
$cell.Style = "Normal"
$ws.Range("E44").Value = "  -0.77%  "

$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = "2.19"
$cell.Style = "Normal"
$ws.Range("E45").Value = "  -2.49%  "

$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = "2.113.40"
$cell.Style = "Normal"
$ws.Range("E46").Value = "  +2.14%  "

$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = "3.27"
$cell.Style = "Normal"
$ws.Range("E47").Value = "  +0.84%  "

$ws.Range("E48").Value = "  +6.83%  "

$cell = $ws.Range("D49")
$cell.NumberFormat = "@"
$cell.Value = "5.43"
$cell.Style = "Normal"
$ws.Range("E49").Value = "  -4.70%  "

$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = "0.904"
$cell.Style = "Normal"
$ws.Range("E50").Value = "  -2.63%  "

$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = "1.33"
$cell.Style = "Normal"
$ws.Range("E51").Value = "  +8.74%  "
